$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 1220.1
$ws.Range("I33").Value = 1635.5714
$ws.Range("J33").Value = 250.66667
$ws.Range("K33").Value = 1635.5714
$ws.Range("L33").Value = 250.66667
$ws.Range("M33").Value = -1406.5714
$ws.Range("N33").Value = -708.6666700000001
$ws.Range("H107").Value = 402.25
$ws.Range("I107").Value = 175.27272
$ws.Range("J107").Value = 901.6
$ws.Range("K107").Value = 175.27272
$ws.Range("L107").Value = 901.6
$ws.Range("M107").Value = 1744.72728
$ws.Range("N107").Value = -4741.6
$ws.Range("H111").Value = 13448.75
$ws.Range("I111").Value = 9318
$ws.Range("J111").Value = 20333.334
$ws.Range("K111").Value = 27954
$ws.Range("L111").Value = 61000.00199999999
$ws.Range("M111").Value = -24887
$ws.Range("N111").Value = -67134.00199999999
$ws.Range("H113").Value = 4432.25
$ws.Range("I113").Value = 3196.6667
$ws.Range("J113").Value = 4961.7856
$ws.Range("K113").Value = 3196.6667
$ws.Range("L113").Value = 4961.7856
$ws.Range("M113").Value = 57.33329999999978
$ws.Range("N113").Value = -11469.7856
$ws.Range("H116").Value = 4999.5835
$ws.Range("I116").Value = 2512.125
$ws.Range("J116").Value = 9974.5
$ws.Range("K116").Value = 2512.125
$ws.Range("L116").Value = 9974.5
$ws.Range("M116").Value = 929.875
$ws.Range("N116").Value = -16858.5
$ws.Range("H138").Value = 1419.2122
$ws.Range("I138").Value = 1074.96
$ws.Range("K138").Value = 3224.88
$ws.Range("M138").Value = 1915.12

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3508.2144
$ws.Range("I2").Value = 2344.6667
$ws.Range("K2").Value = 2344.6667
$ws.Range("M2").Value = -2231.6667
$ws.Range("H32").Value = 6392.8306
$ws.Range("I32").Value = 7038.6313
$ws.Range("J32").Value = 1791.5
$ws.Range("K32").Value = 7038.6313
$ws.Range("L32").Value = 1791.5
$ws.Range("M32").Value = -6751.6313
$ws.Range("N32").Value = -2365.5
$ws.Range("H37").Value = 21994.5
$ws.Range("I37").Value = 19000
$ws.Range("J37").Value = 24989
$ws.Range("K37").Value = 19000
$ws.Range("L37").Value = 24989
$ws.Range("M37").Value = -18727
$ws.Range("N37").Value = -25535
$ws.Range("H116").Value = 3508.2144
$ws.Range("I116").Value = 2344.6667
$ws.Range("K116").Value = 2344.6667
$ws.Range("M116").Value = -50.66670000000022
$ws.Range("H132").Value = 3344.8853
$ws.Range("I132").Value = 1241.8462
$ws.Range("J132").Value = 7073
$ws.Range("K132").Value = 3725.5386
$ws.Range("L132").Value = 21219
$ws.Range("M132").Value = -1195.5386
$ws.Range("N132").Value = -26279

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3508.2144
$ws.Range("I3").Value = 2344.6667
$ws.Range("K3").Value = 2344.6667
$ws.Range("M3").Value = -2230.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2350.6875
$ws.Range("I31").Value = 1738.5312
$ws.Range("J31").Value = 3575
$ws.Range("K31").Value = 1738.5312
$ws.Range("L31").Value = 3575
$ws.Range("M31").Value = -1443.5312
$ws.Range("N31").Value = -4165
$ws.Range("H34").Value = 2350.6875
$ws.Range("I34").Value = 1738.5312
$ws.Range("J34").Value = 3575
$ws.Range("K34").Value = 1738.5312
$ws.Range("L34").Value = 3575
$ws.Range("M34").Value = -1536.5312
$ws.Range("N34").Value = -3979

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 910.6667
$ws.Range("I68").Value = 757.7
$ws.Range("J68").Value = 987.15
$ws.Range("K68").Value = 2273.1
$ws.Range("L68").Value = 2961.45
$ws.Range("M68").Value = -1462.1
$ws.Range("N68").Value = -4583.45
$ws.Range("H71").Value = 910.6667
$ws.Range("I71").Value = 757.7
$ws.Range("J71").Value = 987.15
$ws.Range("K71").Value = 6819.3
$ws.Range("L71").Value = 8884.35
$ws.Range("M71").Value = -2763.3
$ws.Range("N71").Value = -16996.35

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1613.15
$ws.Range("I113").Value = 1622.9412
$ws.Range("J113").Value = 1557.6666
$ws.Range("K113").Value = 1622.9412
$ws.Range("L113").Value = 1557.6666
$ws.Range("M113").Value = 547.0588
$ws.Range("N113").Value = -5897.6666
$ws.Range("H122").Value = 1432818.6
$ws.Range("I122").Value = 5001750
$ws.Range("J122").Value = 5246
$ws.Range("K122").Value = 15005250
$ws.Range("L122").Value = 15738
$ws.Range("M122").Value = -15002800
$ws.Range("N122").Value = -20638

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 70007
$ws.Range("I13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("H61").Value = 3620.8262
$ws.Range("I61").Value = 2481.1177
$ws.Range("J61").Value = 6850
$ws.Range("K61").Value = 2481.1177
$ws.Range("L61").Value = 6850
$ws.Range("M61").Value = -2279.1177
$ws.Range("N61").Value = -7254
$ws.Range("H68").Value = 2757.842
$ws.Range("I68").Value = 2319.8572
$ws.Range("J68").Value = 3013.3333
$ws.Range("K68").Value = 2319.8572
$ws.Range("L68").Value = 3013.3333
$ws.Range("M68").Value = -1570.8572
$ws.Range("N68").Value = -4511.3333
$ws.Range("H71").Value = 2757.842
$ws.Range("I71").Value = 2319.8572
$ws.Range("J71").Value = 3013.3333
$ws.Range("K71").Value = 11599.286
$ws.Range("L71").Value = 15066.6665
$ws.Range("M71").Value = -7855.286
$ws.Range("N71").Value = -22554.6665
$ws.Range("H113").Value = 3620.8262
$ws.Range("I113").Value = 2481.1177
$ws.Range("J113").Value = 6850
$ws.Range("K113").Value = 2481.1177
$ws.Range("L113").Value = 6850
$ws.Range("M113").Value = -311.1176999999998
$ws.Range("N113").Value = -11190
$ws.Range("M13").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 831.6667
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 831.6667
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 831.6667
$ws.Range("N3").Value = -1059.6667
$ws.Range("H107").Value = 659.5789
$ws.Range("I107").Value = 423.7143
$ws.Range("J107").Value = 1320
$ws.Range("K107").Value = 1271.1429
$ws.Range("L107").Value = 3960
$ws.Range("M107").Value = 648.8571000000002
$ws.Range("N107").Value = -7800
$ws.Range("M3").ClearContents()

Write-Output "All edits applied"
